$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column G.
# Both row 3 (643edb51 item) and row 5 (7f4dbbee item) shared the same
# generated-date text, so both must be updated to stay on the same shared string.
$overview.Range("G3").Value = "2016-08-31 14:19:23"
$overview.Range("G5").Value = "2016-08-31 14:19:23"

# zh-cn sheet: rows 3 and 5 shared "ht" / handoff / handback timestamps.
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("E5").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-31 14:19:18"
$zhcn.Range("H5").Value = "2016-08-31 14:19:18"
$zhcn.Range("K3").Value = "2016-08-31 14:19:41"
$zhcn.Range("K5").Value = "2016-08-31 14:19:41"

# de-de sheet: rows 3 and 5 shared "ht" / generate date / handback timestamps.
$dede.Range("E3").Value = "mt"
$dede.Range("E5").Value = "mt"
$dede.Range("H3").Value = "2016-08-31 14:19:23"
$dede.Range("H5").Value = "2016-08-31 14:19:23"
$dede.Range("K3").Value = "2016-08-31 14:19:48"
$dede.Range("K5").Value = "2016-08-31 14:19:48"
